# Apply the "Date and Time" / "Cycle_count" update to the single worksheet.
#
# The new report prepends a "Date and Time" row and inserts a new
# "Cycle Count of battery" row further down; every pre-existing row shifts
# down by one to make room (Excel's row-insert keeps each row's original
# formatting, e.g. the [hh]:mm:ss style on "Total time taken for the ride").
# A number of labels also pick up unit suffixes (e.g. "(kW)", "(C)", "(V)")
# and several metrics were recomputed with the new source data, so each
# row is rewritten to match the target report exactly. Rows whose label
# and value already land correctly after the shift (no textual/numeric
# change vs. the source row) are left untouched to avoid perturbing their
# formatting (e.g. re-triggering autofit row height on the multi-line
# "Mode" cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 1; this shifts the existing rows 1-42 down to 2-43.
$ws.Rows.Item(1).Insert()

$ws.Range('A1').Value = 'Date and Time'
$ws.Range('B1').Value = '2024-03-12 09:19:48.749000 to 2024-03-12 10:20:02.465000'

$ws.Range('A7').Value = 'Starting SoC (%)'
$ws.Range('B7').Value = 96

$ws.Range('A8').Value = 'Ending SoC (%)'
$ws.Range('B8').Value = 9

$ws.Range('A9').Value = 'Total distance covered (km)'
$ws.Range('B9').Value = 38.28985876739389

$ws.Range('A10').Value = 'Total energy consumption(WH/KM)'
$ws.Range('B10').Value = 45.23441700865151

$ws.Range('A11').Value = 'Total SOC consumed(%)'
$ws.Range('B11').Value = 87

$ws.Range('A13').Value = 'Peak Power(kW)'
$ws.Range('B13').Value = 6125.566505

$ws.Range('A14').Value = 'Average Power(kW)'
$ws.Range('B14').Value = -1406.240410300406

$ws.Range('A15').Value = 'Total Energy Regenerated(kWh)'
$ws.Range('B15').Value = 112.3995493272222

$ws.Range('A16').Value = 'Regenerative Effectiveness(%)'
$ws.Range('B16').Value = 6.094035577472369

$ws.Range('A17').Value = 'Highest Cell Voltage(V)'
$ws.Range('B17').Value = 3.337

$ws.Range('A18').Value = 'Lowest Cell Voltage(V)'
$ws.Range('B18').Value = 2.921

$ws.Range('A19').Value = 'Difference in Cell Voltage(V)'
$ws.Range('B19').Value = 0.4160000000000004

$ws.Range('A20').Value = 'Minimum Temperature(C)'
$ws.Range('B20').Value = 21

$ws.Range('A21').Value = 'Maximum Temperature(C)'
$ws.Range('B21').Value = 45

$ws.Range('A22').Value = 'Difference in Temperature(C)'
$ws.Range('B22').Value = 24

$ws.Range('A23').Value = 'Maximum Fet Temperature-BMS(C)'
$ws.Range('B23').Value = 72

$ws.Range('A24').Value = 'Maximum Afe Temperature-BMS(C)'
$ws.Range('B24').Value = 66

$ws.Range('A25').Value = 'Maximum PCB Temperature-BMS(C)'
$ws.Range('B25').Value = 64

$ws.Range('A26').Value = 'Maximum MCU Temperature(C)'
$ws.Range('B26').Value = 51

$ws.Range('A27').Value = 'Maximum Motor Temperature(C)'
$ws.Range('B27').Value = 0

$ws.Range('A28').Value = 'Abnormal Motor Temperature Detected(C)'
$ws.Range('B28').Value = 0

$ws.Range('A29').Value = 'highest cell temp(C)'
$ws.Range('B29').Value = 45

$ws.Range('A30').Value = 'lowest cell temp(C)'
$ws.Range('B30').Value = 21

$ws.Range('A31').Value = 'Difference between Highest and Lowest Cell Temperature at 100% SOC(C)'
$ws.Range('B31').Value = 24

$ws.Range('A32').Value = 'Battery Voltage(V)'
$ws.Range('B32').Value = 55

$ws.Range('A33').Value = 'Total energy charged(kWh)'
$ws.Range('B33').Value = 1.889129229166667

$ws.Range('A34').Value = 'Electricity consumption units(kW)'
$ws.Range('B34').Value = [double]"1.177379670659553e-07"

$ws.Range('A35').Value = 'Cycle Count of battery'
$ws.Range('B35').Value = 137

$ws.Range('A36').Value = 'Idling time percentage'
$ws.Range('B36').Value = 13.48716282092948

$ws.Range('A37').Value = 'Time spent in 0-10 km/h'
$ws.Range('B37').Value = 13.26431124936162

$ws.Range('A38').Value = 'Time spent in 10-20 km/h'
$ws.Range('B38').Value = 5.148799851432286

$ws.Range('A40').Value = 'Time spent in 30-40 km/h'
$ws.Range('B40').Value = 21.60731695993315

$ws.Range('A41').Value = 'Time spent in 40-50 km/h'
$ws.Range('B41').Value = 14.78016621013046

$ws.Range('A42').Value = 'Time spent in 50-60 km/h'
$ws.Range('B42').Value = 9.814754631134221

$ws.Range('A43').Value = 'Time spent in 60-70 km/h'
$ws.Range('B43').Value = 10.51812990389526

$ws.Range('A44').Value = 'Time spent in 70-80 km/h'
$ws.Range('B44').Value = 2.414225358651748

$ws.Range('A45').Value = 'Time spent in 80-90 km/h'
$ws.Range('B45').Value = 0.03714192859464228
